$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C..H on rows 2 and 3 hold docxtemplater-style placeholders for the
# numeric ticket counters (processing/completed/canceled/deferred/closed/
# new_or_reopened). Add the ":formatN()" formatter to each placeholder so the
# exported value is converted/rendered as a number, and give the cells a
# plain integer ("0") number format to match.
$cols = @("C", "D", "E", "F", "G", "H")
foreach ($row in @(2, 3)) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $current = [string]$cell.Value2
        $updated = $current -replace '\}$', ':formatN()}'
        $cell.Value2 = $updated
        $cell.NumberFormat = "0"
    }
}
